$d = $word.ActiveDocument

# Locate the (currently empty) paragraph that sits right after the
# "c) Total de los 100 elementos..." paragraph and right before the
# "d) Copia los 11 elementos..." paragraph.
$paraTotal = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t.Contains("c) Total de los 100")) {
        $paraTotal = $i
    }
}
$paraEmpty = $paraTotal + 1

# 1) Remove the "_GoBack" bookmark currently sitting at the end of the
#    "array = [8] * 5" paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# 2) Replace the (now empty) paragraph that sits between "c) Total..." and
#    "d) Copia..." with the two new paragraphs describing the for-loop, and
#    re-create the "_GoBack" bookmark at the end of the new content.
$target = $d.Paragraphs.Item($paraEmpty)

$xml = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="708"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
    <w:t>for</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
    <w:t xml:space="preserve"> i in </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
    <w:t>arrayc</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:ind w:firstLine="708"/>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>
      <w:iCs/>
    </w:rPr>
    <w:tab/>
    <w:t>suma+=i</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$target.Range.InsertXML($xml)
